$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "66.836.21"
$ws.Range("E2").Value = "  -5.29%  "

# Row 3
$ws.Range("D3").Value = "3.370.68"
$ws.Range("E3").Value = "  -6.76%  "

# Row 4
$ws.Range("E4").Value = "  +0.20%  "

# Row 5
$ws.Range("D5").Value = "'561.42"
$ws.Range("E5").Value = "  -5.99%  "

# Row 6
$ws.Range("D6").Value = "'184.17"
$ws.Range("E6").Value = "  -9.20%  "

# Row 7
$ws.Range("D7").Value = "'0.595"
$ws.Range("E7").Value = "  -5.26%  "

# Row 8
$ws.Range("D8").Value = "'0.999"
$ws.Range("E8").Value = "  -0.01%  "

# Row 9
$ws.Range("D9").Value = "3.364.42"
$ws.Range("E9").Value = "  -6.59%  "

# Row 10
$ws.Range("E10").Value = "  -12.58%  "

# Row 11
$ws.Range("E11").Value = "  -7.55%  "

# Row 12
$ws.Range("D12").Value = "'47.91"
$ws.Range("E12").Value = "  -11.27%  "

# Row 13
$ws.Range("E13").Value = "  -11.44%  "

# Row 14
$ws.Range("E14").Value = "  -9.63%  "

# Row 15
$ws.Range("D15").Value = "3.902.02"
$ws.Range("E15").Value = "  -6.41%  "

# Row 16
$ws.Range("D16").Value = "'610.82"
$ws.Range("E16").Value = "  -9.94%  "

# Row 17
$ws.Range("D17").Value = "66.643.00"
$ws.Range("E17").Value = "  -5.59%  "

# Row 18
$ws.Range("E18").Value = "  -3.75%  "

# Row 19
$ws.Range("D19").Value = "3.367.58"
$ws.Range("E19").Value = "  -7.31%  "

# Row 20
$ws.Range("D20").Value = "'17.68"
$ws.Range("E20").Value = "  -7.87%  "

# Row 21
$ws.Range("D21").Value = "'11.70"
$ws.Range("E21").Value = "  -8.43%  "

# Row 22
$ws.Range("E22").Value = "  -8.72%  "

# Row 23
$ws.Range("D23").Value = "'17.10"
$ws.Range("E23").Value = "  -8.33%  "

# Row 24
$ws.Range("E24").Value = "  -3.90%  "

# Row 25
$ws.Range("D25").Value = "'95.78"
$ws.Range("E25").Value = "  -13.24%  "

# Row 26
$ws.Range("D26").Value = "'4.10"
$ws.Range("E26").Value = "  -10.16%  "

# Row 27
$ws.Range("E27").Value = "  -9.65%  "

# Row 28
$ws.Range("D28").Value = "'9.60"
$ws.Range("E28").Value = "  -9.69%  "

# Row 29
$ws.Range("D29").Value = "'8.83"
$ws.Range("E29").Value = "  -12.93%  "

# Row 30
$ws.Range("E30").Value = "  -10.28%  "

# Row 31
$ws.Range("B31").Value = "dogwifhat"
$ws.Range("C31").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D31").Value = "'3.88"
$ws.Range("E31").Value = "  -13.84%  "

# Row 32
$ws.Range("B32").Value = "NEARProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D32").Value = "'6.37"
$ws.Range("E32").Value = "  -11.25%  "

# Row 33
$ws.Range("D33").Value = "'11.27"
$ws.Range("E33").Value = "  -8.73%  "

# Row 34
$ws.Range("E34").Value = "  -8.32%  "

# Row 35
$ws.Range("D35").Value = "'58.66"
$ws.Range("E35").Value = "  -7.74%  "

# Row 36
$ws.Range("D36").Value = "3.757.20"
$ws.Range("E36").Value = "  -3.35%  "

# Row 37
$ws.Range("D37").Value = "'534.63"
$ws.Range("E37").Value = "  +4.33%  "

# Row 38
$ws.Range("D38").Value = "'0.997"
$ws.Range("E38").Value = "  -0.29%  "

# Row 39
$ws.Range("D39").Value = "'3.75"
$ws.Range("E39").Value = "  +38.38%  "

# Row 40
$ws.Range("D40").Value = "'3.42"
$ws.Range("E40").Value = "  -5.21%  "

# Row 41
$ws.Range("D41").Value = "0.0₃0726"
$ws.Range("E41").Value = "  -15.17%  "

# Row 42
$ws.Range("D42").Value = "'2.73"
$ws.Range("E42").Value = "  -9.81%  "

# Row 43
$ws.Range("D43").Value = "'0.353"
$ws.Range("E43").Value = "  -8.65%  "

# Row 44
$ws.Range("D44").Value = "'0.128"
$ws.Range("E44").Value = "  -7.65%  "

# Row 45
$ws.Range("D45").Value = "'32.76"
$ws.Range("E45").Value = "  -11.05%  "

# Row 46
$ws.Range("D46").Value = "'0.0419"
$ws.Range("E46").Value = "  -10.56%  "

# Row 47
$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").Value = "'3.19"
$ws.Range("E47").Value = "  -7.22%  "

# Row 48
$ws.Range("B48").Value = "ThetaToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D48").Value = "'2.69"
$ws.Range("E48").Value = "  -12.76%  "

# Row 49
$ws.Range("E49").Value = "  -7.92%  "

# Row 50
$ws.Range("D50").Value = "'0.999"
$ws.Range("E50").Value = "  -0.31%  "

# Row 51
$ws.Range("D51").Value = "'7.77"
$ws.Range("E51").Value = "  -9.96%  "
